$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6-9 (no longer needed)
$ws.Range("A6:B9").EntireRow.Delete()

# Update remaining data rows with new control point values
$ws.Range("A2").Value = 21
$ws.Range("B2").Value = 208

$ws.Range("A3").Value = 11
$ws.Range("B3").Value = 130

$ws.Range("A4").Value = 12
$ws.Range("B4").Value = 97

$ws.Range("A5").Value = 22
$ws.Range("B5").Value = 21
